$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by duplicating the "2022-Q2" sheet
#    (same layout/styles) right after "总计", then rename it.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$q2src = $wb.Worksheets.Item("2022-Q2")
$q2src.Copy($null, $total)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Columns D:G store numeric-looking figures as plain TEXT in this
# workbook (matches the source data), so force Text format before
# assigning or Excel will silently coerce them to numbers and lose the
# exact string (e.g. "1.4820" -> 1.482).
$q4.Range("D2:G3").NumberFormat = "@"

# Update the fund data on the new "2022-Q4" sheet with the Q4 figures.
$q4.Cells.Item(2,3).Value = "天弘恒生科技指数（QDII）A"
$q4.Cells.Item(2,4).Value = "39.65"
$q4.Cells.Item(2,5).Value = "93.67"
$q4.Cells.Item(2,6).Value = "3.95"
$q4.Cells.Item(2,7).Value = "1.5662"
$q4.Cells.Item(2,8).Value = 10

$q4.Cells.Item(3,3).Value = "天弘恒生科技指数（QDII）C"
$q4.Cells.Item(3,4).Value = "37.52"
$q4.Cells.Item(3,5).Value = "93.67"
$q4.Cells.Item(3,6).Value = "3.95"
$q4.Cells.Item(3,7).Value = "1.4820"
$q4.Cells.Item(3,8).Value = 10

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: add a new data row for 2022-Q4 and
#    push the existing quarters down by one row (2020-Q4 becomes row 9).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")

# Grow the table by one row, copying formats from the current last row
# (row 8) down into the new row 9 so style/number formatting match.
$ws.Range("A8:D8").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)

# Re-write all the data rows (2..9) bottom-up so every quarter/metric
# lines up with its final row number.
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "2020-Q4"
$ws.Cells.Item(9,3).Value = 5
$ws.Cells.Item(9,4).Value = 0.41

$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "2021-Q1"
$ws.Cells.Item(8,3).Value = 5
$ws.Cells.Item(8,4).Value = 0.64

$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "2021-Q2"
$ws.Cells.Item(7,3).Value = 7
$ws.Cells.Item(7,4).Value = 1.25

$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "2021-Q3"
$ws.Cells.Item(6,3).Value = 9
$ws.Cells.Item(6,4).Value = 2.26

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "2021-Q4"
$ws.Cells.Item(5,3).Value = 7
$ws.Cells.Item(5,4).Value = 2.58

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "2022-Q1"
$ws.Cells.Item(4,3).Value = 4
$ws.Cells.Item(4,4).Value = 3.07

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "2022-Q2"
$ws.Cells.Item(3,3).Value = 2
$ws.Cells.Item(3,4).Value = 2.87

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "2022-Q4"
$ws.Cells.Item(2,3).Value = 2
$ws.Cells.Item(2,4).Value = 3.05
